$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2022" column (N) to the table, mirroring the formatting already
# used by the "2021" column (M) for the header and each data row.
$ws.Range("M4:M13").Copy($ws.Range("N4"))

$ws.Range("N4").Value = 2022
$ws.Range("N5").Value = 4.3
$ws.Range("N6").Value = 5.0999999999999996
$ws.Range("N7").Value = 3.1
$ws.Range("N8").Value = 2.9
$ws.Range("N9").Value = 3.4
$ws.Range("N10").Value = 2.2999999999999998
$ws.Range("N11").Value = 92.8
$ws.Range("N12").Value = 91.6
$ws.Range("N13").Value = 94.6

$ws.Range("N15").Select()
